# Auto-generated edit script: updates Leve market-price/profit figures
# across the per-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3416.5
$ws.Range("I18").Value = 3500
$ws.Range("K18").Value = 3500
$ws.Range("M18").Value = -3216

$ws.Range("H32").Value = 4539.857
$ws.Range("J32").Value = 6444.5
$ws.Range("L32").Value = 6444.5
$ws.Range("N32").Value = -7096.5

$ws.Range("H70").Value = 1224.25
$ws.Range("J70").Value = 1369.9
$ws.Range("L70").Value = 4109.700000000001
$ws.Range("N70").Value = -4649.700000000001

$ws.Range("H73").Value = 1224.25
$ws.Range("J73").Value = 1369.9
$ws.Range("L73").Value = 4109.700000000001
$ws.Range("N73").Value = -5981.700000000001

$ws.Range("H99").Value = 599.2
$ws.Range("I99").Value = 278.5
$ws.Range("J99").Value = 813
$ws.Range("K99").Value = 835.5
$ws.Range("L99").Value = 2439
$ws.Range("M99").Value = 662.5
$ws.Range("N99").Value = -5435

$ws.Range("H101").Value = 25000374
$ws.Range("I101").Value = 33333500
$ws.Range("J101").Value = 999
$ws.Range("K101").Value = 100000500
$ws.Range("L101").Value = 2997
$ws.Range("M101").Value = -99998878
$ws.Range("N101").Value = -6241

$ws.Range("H113").Value = 3339
$ws.Range("J113").Value = 2953
$ws.Range("L113").Value = 2953
$ws.Range("N113").Value = -9461

$ws.Range("H116").Value = 4666.6665
$ws.Range("I116").Value = 4666.6665
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 4666.6665
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1224.6665
$ws.Range("N116").Value = ""

$ws.Range("H132").Value = 1640.9642
$ws.Range("I132").Value = 1331.4073
$ws.Range("K132").Value = 3994.2219
$ws.Range("M132").Value = -1464.2219

$ws.Range("H138").Value = 3267.9285
$ws.Range("I138").Value = 3266.3845
$ws.Range("K138").Value = 9799.1535
$ws.Range("M138").Value = -4659.1535


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3068.7778
$ws.Range("I32").Value = 3156.577
$ws.Range("K32").Value = 3156.577
$ws.Range("M32").Value = -2869.577

$ws.Range("H122").Value = 1521.8096
$ws.Range("I122").Value = 868.1429000000001
$ws.Range("K122").Value = 2604.4287
$ws.Range("M122").Value = -154.4287000000004


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 19900
$ws.Range("J97").Value = 19900
$ws.Range("L97").Value = 19900
$ws.Range("N97").Value = -21882

$ws.Range("H99").Value = 1005.3333
$ws.Range("I99").Value = 942.63635
$ws.Range("J99").Value = 1695
$ws.Range("K99").Value = 942.63635
$ws.Range("L99").Value = 1695
$ws.Range("M99").Value = 555.36365
$ws.Range("N99").Value = -4691

$ws.Range("H130").Value = 83999.2
$ws.Range("J130").Value = 83999.2
$ws.Range("L130").Value = 83999.2
$ws.Range("N130").Value = -94039.2

$ws.Range("H134").Value = 2666.0435
$ws.Range("I134").Value = 2719.1365
$ws.Range("J134").Value = 1498
$ws.Range("K134").Value = 8157.4095
$ws.Range("L134").Value = 4494
$ws.Range("M134").Value = -5622.4095
$ws.Range("N134").Value = -9564


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = ""

$ws.Range("H31").Value = 2322.9285
$ws.Range("J31").Value = 1949.8334
$ws.Range("L31").Value = 1949.8334
$ws.Range("N31").Value = -2539.8334

$ws.Range("H34").Value = 2322.9285
$ws.Range("J34").Value = 1949.8334
$ws.Range("L34").Value = 1949.8334
$ws.Range("N34").Value = -2353.8334

$ws.Range("H62").Value = 5050
$ws.Range("J62").Value = 5250
$ws.Range("L62").Value = 5250
$ws.Range("N62").Value = -6498

$ws.Range("H65").Value = 5050
$ws.Range("J65").Value = 5250
$ws.Range("L65").Value = 26250
$ws.Range("N65").Value = -32490

$ws.Range("H94").Value = 93151.62
$ws.Range("I94").Value = 164566.28
$ws.Range("J94").Value = 9834.5
$ws.Range("K94").Value = 164566.28
$ws.Range("L94").Value = 9834.5
$ws.Range("M94").Value = -164115.28
$ws.Range("N94").Value = -10736.5

$ws.Range("H99").Value = 1366.3334
$ws.Range("I99").Value = 1366.3334
$ws.Range("K99").Value = 1366.3334
$ws.Range("M99").Value = 131.6666

$ws.Range("H126").Value = 1366.3334
$ws.Range("I126").Value = 1366.3334
$ws.Range("K126").Value = 4099.0002
$ws.Range("M126").Value = -1629.0002

$ws.Range("H134").Value = 919.95
$ws.Range("I134").Value = 919.95
$ws.Range("K134").Value = 2759.85
$ws.Range("M134").Value = -224.8500000000004


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 726.3333
$ws.Range("I92").Value = 454.33334
$ws.Range("J92").Value = 998.3333
$ws.Range("K92").Value = 1363.00002
$ws.Range("L92").Value = 2994.9999
$ws.Range("M92").Value = -115.0000199999999
$ws.Range("N92").Value = -5490.9999

$ws.Range("H106").Value = 19092.334
$ws.Range("I106").Value = 15555
$ws.Range("K106").Value = 46665
$ws.Range("M106").Value = -45719

$ws.Range("H128").Value = 646389.2
$ws.Range("I128").Value = 646389.2
$ws.Range("K128").Value = 1939167.6
$ws.Range("M128").Value = -1934187.6

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = ""


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3179.625
$ws.Range("I122").Value = 2970.6365
$ws.Range("K122").Value = 8911.9095
$ws.Range("M122").Value = -6461.9095

$ws.Range("H126").Value = 9000
$ws.Range("I126").Value = 9000
$ws.Range("K126").Value = 27000
$ws.Range("M126").Value = -24530

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 101856
$ws.Range("J46").Value = 1800
$ws.Range("L46").Value = 1800
$ws.Range("N46").Value = -2176


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""

$ws.Range("H62").Value = 13350.375
$ws.Range("J62").Value = 9833.5
$ws.Range("L62").Value = 9833.5
$ws.Range("N62").Value = -11081.5

$ws.Range("H65").Value = 13350.375
$ws.Range("J65").Value = 9833.5
$ws.Range("L65").Value = 49167.5
$ws.Range("N65").Value = -55407.5

$ws.Range("H107").Value = 933.8125
$ws.Range("I107").Value = 894.8
$ws.Range("J107").Value = 998.8333
$ws.Range("K107").Value = 2684.4
$ws.Range("L107").Value = 2996.4999
$ws.Range("M107").Value = -764.3999999999996
$ws.Range("N107").Value = -6836.4999

$ws.Range("H126").Value = 2192.0667
$ws.Range("I126").Value = 2192.0667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6576.2001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4106.2001
$ws.Range("N126").Value = ""

$ws.Range("H132").Value = 5864.4
$ws.Range("I132").Value = 6113.3076
$ws.Range("J132").Value = 4246.5
$ws.Range("K132").Value = 18339.9228
$ws.Range("L132").Value = 12739.5
$ws.Range("M132").Value = -15809.9228
$ws.Range("N132").Value = -17799.5

$ws.Range("H136").Value = 10017.8125
$ws.Range("I136").Value = 7352
$ws.Range("J136").Value = 50005
$ws.Range("K136").Value = 22056
$ws.Range("L136").Value = 150015
$ws.Range("M136").Value = -19506
$ws.Range("N136").Value = -155115

